$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber the id sequence in column A for rows 44:50 ---
$ws.Range("A44").Value = 4
$ws.Range("A45").Value = 5
$ws.Range("A46").Value = 6
$ws.Range("A47").Value = 7
$ws.Range("A48").Value = 8
$ws.Range("A49").Value = 9
$ws.Range("A50").Value = 10

# --- Drop the now-unused tail rows (51:58): wipe A/B, blank out C/D ---
$ws.Range("A51:D58").ClearContents()

# --- Hyperlinks tied to the blanked D51:D58 cells must go too. This host's
# Range(...).Hyperlinks.Delete() removes every hyperlink on the sheet
# regardless of the range it was called on, so capture the addresses that
# still need to exist (D41:D50), delete everything, then recreate just
# those - putting each cell's style back to the shared "Hyperlink" xf
# (plus vertical-center where that was already the case) so we don't
# leave a stray one-off format behind. ---
$addrs = @{
    41 = "mailto:abc@123"; 42 = "mailto:abc@432"; 43 = "mailto:fds@345";
    44 = "mailto:abc@123"; 45 = "mailto:abc@432"; 46 = "mailto:fds@345";
    47 = "mailto:abc@123"; 48 = "mailto:abc@432"; 49 = "mailto:fds@345";
    50 = "mailto:abc@123"
}
$centered = @(42, 45, 48)

$ws.Range("A1").Hyperlinks.Delete()

foreach ($r in 41..50) {
    $cell = $ws.Cells.Item($r, 4)
    $ws.Hyperlinks.Add($cell, $addrs[$r])
    $cell.Style = "Hyperlink"
    if ($centered -contains $r) {
        $cell.VerticalAlignment = -4108
    }
}

# --- Sheet view cursor moved off the old B59:B60 selection ---
$ws.Range("O57").Select() | Out-Null
